$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change Promo.promo_code Domain from "NA" to "UNIQUE" (row 29, column M)
$ws.Cells.Item(29, 13).Value = "UNIQUE"

# New "Customer" table block, rows 43-46
# Row 43: id (PK)
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = "Customer"
$ws.Cells.Item(43, 3).Value = "id"
$ws.Cells.Item(43, 4).Value = "NA"
$ws.Cells.Item(43, 5).Value = "NA"
$ws.Cells.Item(43, 6).Value = "int"
$ws.Cells.Item(43, 7).Value = $true
$ws.Cells.Item(43, 8).Value = $false
$ws.Cells.Item(43, 9).Value = $true
$ws.Cells.Item(43, 10).Value = $false
$ws.Cells.Item(43, 11).Value = $false
$ws.Cells.Item(43, 12).Value = $false
$ws.Cells.Item(43, 13).Value = "UNIQUE"
$ws.Cells.Item(43, 14).Value = "Standard Auto-Increment PK"
$ws.Cells.Item(43, 15).Value = "Name and email for a customer who will recieve promo emails"

# Row 44: first_name
$ws.Cells.Item(44, 1).Value = " "
$ws.Cells.Item(44, 2).Value = " "
$ws.Cells.Item(44, 3).Value = "first_name"
$ws.Cells.Item(44, 4).Value = "NA"
$ws.Cells.Item(44, 5).Value = 200
$ws.Cells.Item(44, 6).Value = "nvarchar"
$ws.Cells.Item(44, 7).Value = $false
$ws.Cells.Item(44, 8).Value = $false
$ws.Cells.Item(44, 9).Value = $true
$ws.Cells.Item(44, 10).Value = $false
$ws.Cells.Item(44, 11).Value = $false
$ws.Cells.Item(44, 12).Value = $false
$ws.Cells.Item(44, 13).Value = "NA"
$ws.Cells.Item(44, 15).Value = " "

# Row 45: last_name
$ws.Cells.Item(45, 1).Value = " "
$ws.Cells.Item(45, 2).Value = " "
$ws.Cells.Item(45, 3).Value = "last_name"
$ws.Cells.Item(45, 4).Value = "NA"
$ws.Cells.Item(45, 5).Value = 200
$ws.Cells.Item(45, 6).Value = "nvarchar"
$ws.Cells.Item(45, 7).Value = $false
$ws.Cells.Item(45, 8).Value = $false
$ws.Cells.Item(45, 9).Value = $true
$ws.Cells.Item(45, 10).Value = $false
$ws.Cells.Item(45, 11).Value = $false
$ws.Cells.Item(45, 12).Value = $false
$ws.Cells.Item(45, 13).Value = "NA"
$ws.Cells.Item(45, 15).Value = " "

# Row 46: email
$ws.Cells.Item(46, 1).Value = " "
$ws.Cells.Item(46, 2).Value = " "
$ws.Cells.Item(46, 3).Value = "email"
$ws.Cells.Item(46, 4).Value = "NA"
$ws.Cells.Item(46, 5).Value = 254
$ws.Cells.Item(46, 6).Value = "nvarchar"
$ws.Cells.Item(46, 7).Value = $false
$ws.Cells.Item(46, 8).Value = $false
$ws.Cells.Item(46, 9).Value = $true
$ws.Cells.Item(46, 10).Value = $false
$ws.Cells.Item(46, 11).Value = $false
$ws.Cells.Item(46, 12).Value = $false
$ws.Cells.Item(46, 13).Value = "NA"
$ws.Cells.Item(46, 15).Value = " "

# Resize the table to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:O48"))

# Widen column O to fit the longer "Table Desc" text
$ws.Columns.Item(15).ColumnWidth = 73.6
